$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1104
$ws.Range("F7").Value = 136
$ws.Range("F8").Value = 1395
$ws.Range("F9").Value = 62
$ws.Range("F10").Value = 87
$ws.Range("F11").Value = 610
$ws.Range("F12").Value = 125
$ws.Range("F13").Value = 79
$ws.Range("F14").Value = 1300
$ws.Range("F15").Value = 443
$ws.Range("F17").Value = 125
$ws.Range("F19").Value = 670
$ws.Range("F20").Value = 2531
$ws.Range("F22").Value = 42
$ws.Range("F25").Value = 279
$ws.Range("F27").Value = 7
$ws.Range("F29").Value = 559
$ws.Range("F30").Value = 912
$ws.Range("F32").Value = 39
$ws.Range("F34").Value = 159
$ws.Range("F35").Value = 25
$ws.Range("F36").Value = 234

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 712
$ws.Range("F5").Value = 600
$ws.Range("F6").Value = 600
$ws.Range("F15").Value = 308
$ws.Range("F16").Value = 308
$ws.Range("F17").Value = 64
$ws.Range("F22").Value = 596
$ws.Range("F24").Value = 22
$ws.Range("F26").Value = 217
$ws.Range("F27").Value = 218
$ws.Range("F31").Value = 16

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2042
$ws.Range("F6").Value = 2256
$ws.Range("F10").Value = 1090
$ws.Range("F11").Value = 227

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2042
$ws.Range("F4").Value = 2256
$ws.Range("F10").Value = 1090
$ws.Range("F12").Value = 227
$ws.Range("F14").Value = 712
$ws.Range("F15").Value = 1104
$ws.Range("F17").Value = 136
$ws.Range("F18").Value = 1395
$ws.Range("F19").Value = 600
$ws.Range("F20").Value = 62
$ws.Range("F21").Value = 87
$ws.Range("F22").Value = 610
$ws.Range("F23").Value = 125
$ws.Range("F25").Value = 79
$ws.Range("F26").Value = 1300
$ws.Range("F27").Value = 443
$ws.Range("F30").Value = 670
$ws.Range("F31").Value = 2531
$ws.Range("F34").Value = 279
$ws.Range("F39").Value = 559
$ws.Range("F40").Value = 912
$ws.Range("F41").Value = 308
$ws.Range("F43").Value = 39
$ws.Range("F46").Value = 22
$ws.Range("F47").Value = 217
$ws.Range("F49").Value = 159
$ws.Range("F50").Value = 25
$ws.Range("F51").Value = 234
